$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1:E6").Value = "01.01.1901"
$ws.Range("E1").Select() | Out-Null
